# Se cambio prediccion minima a 0.5
# Update the prediction probability values for the affected rows and
# re-color the cells that now cross the (lowered) minimum prediction
# threshold from "reprobado" (red) to "aprobado" (green).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$greenColor = 6154589   # BGR for fill FF5DE95D (Aprobado / pass)

# Rows whose probability value changed AND whose fill must flip from
# red (below old 0.55 threshold) to green (now above new 0.5 threshold)
$ws.Range("G14").Value = 0.548
$ws.Range("G14").Interior.Color = $greenColor

$ws.Range("G26").Value = 0.548
$ws.Range("G26").Interior.Color = $greenColor

$ws.Range("G42").Value = 0.548
$ws.Range("G42").Interior.Color = $greenColor

$ws.Range("G44").Value = 0.548
$ws.Range("G44").Interior.Color = $greenColor

# Rows whose probability value changed but were already green (stay green)
$ws.Range("G16").Value = 0.928
$ws.Range("G23").Value = 0.888
$ws.Range("G27").Value = 0.708
$ws.Range("G34").Value = 0.708
$ws.Range("G47").Value = 0.972

# Summary "Aprobados" (approved) count increased from 45 to 49
$ws.Range("G52").Value = 49
